# Natmi following Dr Hou advice
#
# Bmp7-Bmpr2 LR-pairs sheet: recompute the 3-row FAPs-source block with the
# refreshed specificity / expression numbers, and add the new sCs-source
# block (rows 5-7) so FAPs and sCs are both represented as sending clusters
# against the ECs / FAPs / sCs target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the various detection-rate / expression / specificity numbers.
$rows = @(
  @(2,  "FAPs", "Bmp7", "Bmpr2", "ECs",  3, 1, 1.668521,            5.005563,    0.9677024783929865, 0.9677024783929865, 3, 1, 40.70766766666667,  122.123003,         0.3776398983502007, 0.3776398983502007, 67.92159836285434,  611.294385265689,   0.3654430655735648,  0.3654430655735647),
  @(3,  "FAPs", "Bmp7", "Bmpr2", "FAPs", 3, 1, 1.668521,            5.005563,    0.9677024783929865, 0.9677024783929865, 3, 1, 39.715023,          119.145069,         0.3684312589831062, 0.3684312589831062, 66.265349890983,    596.3881490188471,  0.3565318424354001,  0.3565318424354001),
  @(4,  "FAPs", "Bmp7", "Bmpr2", "sCs",  3, 1, 1.668521,            5.005563,    0.9677024783929865, 0.9677024783929865, 3, 1, 27.37224266666666,  82.11672799999999, 0.253928842666693,  0.253928842666693,  45.67116170642934,  411.040455357864,   0.2457275703840216,  0.2457275703840216),
  @(5,  "sCs",  "Bmp7", "Bmpr2", "ECs",  1, 0.3333333333333333, 0.05568766666666666, 0.167063, 0.03229752160701353, 0.03229752160701353, 3, 1, 40.70766766666667,  122.123003,  0.3776398983502007, 0.3776398983502007, 2.266915027798778,  20.402235250189,    0.012196832776636,   0.012196832776636),
  @(6,  "sCs",  "Bmp7", "Bmpr2", "FAPs", 1, 0.3333333333333333, 0.05568766666666666, 0.167063, 0.03229752160701353, 0.03229752160701353, 3, 1, 39.715023,          119.145069,  0.3684312589831062, 0.3684312589831062, 2.211636962483,     19.904732662347,    0.01189941654770607, 0.01189941654770607),
  @(7,  "sCs",  "Bmp7", "Bmpr2", "sCs",  1, 0.3333333333333333, 0.05568766666666666, 0.167063, 0.03229752160701353, 0.03229752160701353, 3, 1, 27.37224266666666,  82.11672799999999, 0.253928842666693,  0.253928842666693,  1.524296325540444,  13.718666929864,    0.008201272282671458, 0.008201272282671458)
)

foreach ($rowSpec in $rows) {
  $rowNum = $rowSpec[0]
  for ($colIdx = 1; $colIdx -lt $rowSpec.Length; $colIdx++) {
    $colNum = $colIdx  # rowSpec[1] -> column A (1), rowSpec[2] -> column B (2), ...
    $ws.Cells.Item($rowNum, $colNum).Value = $rowSpec[$colIdx]
  }
}

Write-Host "Wrote rows 2-7 of Sheet1 (dimension now A1:T7)."
